$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '95.261.09'
$ws.Range("E2").Value = '  -2.10%  '

$ws.Range("D3").Value = '3.593.15'
$ws.Range("E3").Value = '  -3.41%  '

$ws.Range("D4").Value = '2.61'
$ws.Range("E4").Value = '  +35.74%  '

$ws.Range("E5").Value = '  +0.21%  '

$ws.Range("D6").Value = '221.45'
$ws.Range("E6").Value = '  -6.54%  '

$ws.Range("D7").Value = '632.37'
$ws.Range("E7").Value = '  -4.01%  '

$ws.Range("D8").Value = '0.413'
$ws.Range("E8").Value = '  -5.09%  '

$ws.Range("D9").Value = '1.16'
$ws.Range("E9").Value = '  +8.59%  '

$ws.Range("E10").Value = '  +0.13%  '

$ws.Range("D11").Value = '3.589.22'
$ws.Range("E11").Value = '  -3.46%  '

$ws.Range("D12").Value = '47.24'
$ws.Range("E12").Value = '  +5.39%  '

$ws.Range("E13").Value = '  +0.76%  '

$ws.Range("E14").Value = '  -9.74%  '

$ws.Range("D15").Value = '6.42'
$ws.Range("E15").Value = '  -7.35%  '

$ws.Range("D16").Value = '4.266.00'
$ws.Range("E16").Value = '  -3.27%  '

$ws.Range("D17").Value = '94.832.76'
$ws.Range("E17").Value = '  -2.05%  '

$ws.Range("D18").Value = '21.96'
$ws.Range("E18").Value = '  +17.41%  '

$ws.Range("D19").Value = '8.81'
$ws.Range("E19").Value = '  -2.26%  '

$ws.Range("D20").Value = '13.70'
$ws.Range("E20").Value = '  +5.26%  '

$ws.Range("D21").Value = '3.591.89'
$ws.Range("E21").Value = '  -3.61%  '

$ws.Range("D22").Value = '0.535'
$ws.Range("E22").Value = '  +4.85%  '

$ws.Range("D23").Value = '0.277'
$ws.Range("E23").Value = '  +45.58%  '

$ws.Range("D24").Value = '507.80'
$ws.Range("E24").Value = '  -3.48%  '

$ws.Range("D25").Value = '3.21'
$ws.Range("E25").Value = '  -7.80%  '

$ws.Range("D26").Value = '118.94'
$ws.Range("E26").Value = '  +11.97%  '

$ws.Range("D27").Value = '0.0000197'
$ws.Range("E27").Value = '  -11.84%  '

$ws.Range("D28").Value = '6.74'
$ws.Range("E28").Value = '  -2.28%  '

$ws.Range("D29").Value = '3.779.25'
$ws.Range("E29").Value = '  -3.43%  '

$ws.Range("D30").Value = '12.55'
$ws.Range("E30").Value = '  -7.55%  '

$ws.Range("D31").Value = '12.62'
$ws.Range("E31").Value = '  -0.56%  '

$ws.Range("D32").Value = '2.98'
$ws.Range("E32").Value = '  -1.36%  '

$ws.Range("E33").Value = '  +0.01%  '

$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("B35").Value = 'Cronos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D35").Value = '0.178'
$ws.Range("E35").Value = '  -7.13%  '

$ws.Range("D36").Value = '32.13'
$ws.Range("E36").Value = '  -1.54%  '

$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").Value = '0.606'
$ws.Range("E37").Value = '  +2.37%  '

$ws.Range("D38").Value = '1.73'
$ws.Range("E38").Value = '  -6.01%  '

$ws.Range("E39").Value = '  +0.00%  '

$ws.Range("D40").Value = '8.25'
$ws.Range("E40").Value = '  -5.95%  '

$ws.Range("D41").Value = '571.77'
$ws.Range("E41").Value = '  -10.97%  '

$ws.Range("E42").Value = '  +3.10%  '

$ws.Range("D43").Value = '41.41'
$ws.Range("E43").Value = '  +1.68%  '

$ws.Range("E44").Value = '  +0.53%  '

$ws.Range("D45").Value = '0.0496'
$ws.Range("E45").Value = '  +9.06%  '

$ws.Range("E46").Value = '  -7.59%  '

$ws.Range("D47").Value = '0.944'
$ws.Range("E47").Value = '  -2.48%  '

$ws.Range("D48").Value = '1.91'
$ws.Range("E48").Value = '  -5.82%  '

$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = '8.77'
$ws.Range("E49").Value = '  +1.11%  '

$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '225.05'
$ws.Range("E50").Value = '  +8.93%  '

$ws.Range("D51").Value = '23.51'
$ws.Range("E51").Value = '  -0.57%  '

